# Weekly price-sheet update: a new observation is inserted as the first
# data row (row 12, right after the fixed header block in rows 1-11),
# pushing the existing data rows 12-87 down to 13-88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 12; everything below (old rows 12-87)
# shifts down to rows 13-88.
$ws.Rows.Item(12).Insert()

# Populate the new row with the latest weekly observation.
$ws.Cells.Item(12, 1).Value  = 10
$ws.Cells.Item(12, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(12, 3).Value  = "La Araucanía"
$ws.Cells.Item(12, 4).Value  = 45169
$ws.Cells.Item(12, 5).Value  = 9
$ws.Cells.Item(12, 6).Value  = 300000000
$ws.Cells.Item(12, 7).Value  = "Espárragos"
$ws.Cells.Item(12, 8).Value  = "Sin especificar"
$ws.Cells.Item(12, 9).Value  = "Primera"
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 3000
$ws.Cells.Item(12, 12).Value = 3000
$ws.Cells.Item(12, 13).Value = 3000
$ws.Cells.Item(12, 14).Value = "$/kilo"
$ws.Cells.Item(12, 15).Value = "Región del Maule"
$ws.Cells.Item(12, 16).Value = 3000
$ws.Cells.Item(12, 17).Value = 1
$ws.Cells.Item(12, 18).Value = "Hortaliza"
